$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.499.45"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "2.379.17"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.74%  "
$ws.Range("E7").Value = "  +0.82%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -1.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0922"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.987"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.03%  "
$ws.Range("D15").Value = "2.741.55"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.16%  "
$ws.Range("D17").Value = "2.362.32"
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("D18").Value = "45.442.22"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.31"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +14.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.72%  "
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("E30").Value = "  -1.87%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0953"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "37.56"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "168.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.81"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.77%  "
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("E36").Value = "  -2.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.14%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.94%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0356"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "70.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "97.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.61%  "
$ws.Range("E45").Value = "  -4.34%  "
$ws.Range("D46").Value = "1.867.07"
$ws.Range("E46").Value = "  +13.49%  "
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("E48").Value = "  +1.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "84.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.34%  "
$ws.Range("E51").Value = "  -0.69%  "
